# Apply Betfair odds updates for 2025-11-25 (rows 2-22, columns F:AO).
# Values below are taken from the authoritative before/after cell diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("AD2").Value = 20
$ws.Range("AE2").Value = 60
$ws.Range("AG2").Value = 10
$ws.Range("AJ2").Value = 17.5
$ws.Range("AK2").Value = 15.5
$ws.Range("AM2").Value = 80
$ws.Range("AO2").Value = 50
$ws.Range("F2").Value = 1.7
$ws.Range("G2").Value = 1.72
$ws.Range("H2").Value = 5.1
$ws.Range("I2").Value = 5.3
$ws.Range("P2").Value = 2.46
$ws.Range("Q2").Value = 1.66
$ws.Range("S2").Value = 2.64
$ws.Range("V2").Value = 1.23
$ws.Range("W2").Value = 2.38
$ws.Range("X2").Value = 23
$ws.Range("Z2").Value = 44
# Row 3
$ws.Range("AJ3").Value = 14
$ws.Range("AN3").Value = 9.800000000000001
$ws.Range("F3").Value = 1.4
$ws.Range("H3").Value = 9.800000000000001
$ws.Range("I3").Value = 11.5
$ws.Range("J3").Value = 4.7
$ws.Range("K3").Value = 5.2
$ws.Range("L3").Value = 1.41
$ws.Range("N3").Value = 3.55
$ws.Range("P3").Value = 1.88
$ws.Range("R3").Value = 1.33
$ws.Range("V3").Value = 1.09
$ws.Range("W3").Value = 3.2
# Row 4
$ws.Range("AL4").Value = 50
$ws.Range("AO4").Value = 11
$ws.Range("G4").Value = 4.4
$ws.Range("H4").Value = 1.94
$ws.Range("I4").Value = 1.95
$ws.Range("K4").Value = 4
$ws.Range("L4").Value = 1.33
$ws.Range("N4").Value = 4.6
$ws.Range("S4").Value = 2.92
$ws.Range("V4").Value = 2.04
$ws.Range("W4").Value = 1.3
# Row 5
$ws.Range("AC5").Value = 10
$ws.Range("AF5").Value = 11.5
$ws.Range("AL5").Value = 26
$ws.Range("F5").Value = 1.67
$ws.Range("K5").Value = 4.6
$ws.Range("L5").Value = 1.3
$ws.Range("N5").Value = 5.8
# Row 6
$ws.Range("AD6").Value = 34
$ws.Range("AL6").Value = 27
$ws.Range("AM6").Value = 60
$ws.Range("AO6").Value = 980
$ws.Range("I6").Value = 7.4
$ws.Range("R6").Value = 2.06
$ws.Range("S6").Value = 1.7
$ws.Range("V6").Value = 1.16
$ws.Range("X6").Value = 55
# Row 7
$ws.Range("G7").Value = 2.16
$ws.Range("Q7").Value = 1.48
$ws.Range("W7").Value = 1.86
# Row 8
$ws.Range("G8").Value = 2.14
$ws.Range("U8").Value = 2.54
# Row 9
$ws.Range("AG9").Value = 13
$ws.Range("AJ9").Value = 17
$ws.Range("F9").Value = 1.42
$ws.Range("G9").Value = 1.49
$ws.Range("K9").Value = 6.2
$ws.Range("N9").Value = 6.8
$ws.Range("P9").Value = 2.98
$ws.Range("R9").Value = 1.8
$ws.Range("T9").Value = 1.65
$ws.Range("W9").Value = 3
# Row 10
$ws.Range("Q10").Value = 2.54
# Row 11
$ws.Range("J11").Value = 3.4
$ws.Range("U11").Value = 2.22
$ws.Range("W11").Value = 1.5
# Row 12
$ws.Range("I12").Value = 5.9
$ws.Range("N12").Value = 2.86
$ws.Range("O12").Value = 1.43
$ws.Range("T12").Value = 2.04
$ws.Range("V12").Value = 1.21
# Row 13
$ws.Range("AC13").Value = 8.800000000000001
$ws.Range("AH13").Value = 16.5
$ws.Range("AK13").Value = 17.5
$ws.Range("AM13").Value = 75
$ws.Range("AN13").Value = 9.6
$ws.Range("N13").Value = 4.8
$ws.Range("O13").Value = 1.24
$ws.Range("P13").Value = 2.32
$ws.Range("Q13").Value = 1.73
$ws.Range("R13").Value = 1.52
$ws.Range("S13").Value = 2.84
$ws.Range("T13").Value = 1.69
$ws.Range("U13").Value = 2.36
$ws.Range("X13").Value = 18.5
$ws.Range("Y13").Value = 19.5
$ws.Range("Z13").Value = 36
# Row 14
$ws.Range("AA14").Value = 28
$ws.Range("AO14").Value = 11
$ws.Range("F14").Value = 3.4
$ws.Range("G14").Value = 3.45
$ws.Range("L14").Value = 1.28
$ws.Range("P14").Value = 2.52
$ws.Range("Q14").Value = 1.64
# Row 15
$ws.Range("AC15").Value = 11
$ws.Range("AH15").Value = 12.5
$ws.Range("AJ15").Value = 36
$ws.Range("AN15").Value = 9.199999999999999
$ws.Range("H15").Value = 2.86
$ws.Range("L15").Value = 1.22
$ws.Range("R15").Value = 1.94
$ws.Range("S15").Value = 2.02
$ws.Range("U15").Value = 3.3
$ws.Range("W15").Value = 1.68
$ws.Range("X15").Value = 34
$ws.Range("Y15").Value = 24
# Row 16
$ws.Range("AA16").Value = 570
$ws.Range("AB16").Value = 14.5
$ws.Range("AF16").Value = 9.800000000000001
$ws.Range("AJ16").Value = 10.5
$ws.Range("AN16").Value = 3.15
$ws.Range("J16").Value = 7.6
$ws.Range("P16").Value = 3.4
$ws.Range("R16").Value = 1.98
$ws.Range("S16").Value = 1.96
$ws.Range("T16").Value = 1.88
$ws.Range("X16").Value = 40
# Row 17
$ws.Range("AA17").Value = 550
$ws.Range("AE17").Value = 190
$ws.Range("F17").Value = 1.31
$ws.Range("G17").Value = 1.32
$ws.Range("R17").Value = 1.61
$ws.Range("S17").Value = 2.54
$ws.Range("W17").Value = 4.1
# Row 18
$ws.Range("AN18").Value = 24
$ws.Range("L18").Value = 1.36
$ws.Range("P18").Value = 2.14
$ws.Range("Q18").Value = 1.84
$ws.Range("S18").Value = 3.1
$ws.Range("W18").Value = 1.49
# Row 19
$ws.Range("AJ19").Value = 60
$ws.Range("AK19").Value = 40
$ws.Range("AL19").Value = 50
$ws.Range("AO19").Value = 24
$ws.Range("F19").Value = 3.2
$ws.Range("G19").Value = 3.25
$ws.Range("H19").Value = 2.48
$ws.Range("I19").Value = 2.52
$ws.Range("N19").Value = 3.6
$ws.Range("P19").Value = 1.86
$ws.Range("V19").Value = 1.66
$ws.Range("W19").Value = 1.44
$ws.Range("X19").Value = 13
# Row 20
$ws.Range("F20").Value = 4.5
$ws.Range("H20").Value = 1.93
$ws.Range("I20").Value = 1.98
$ws.Range("J20").Value = 3.55
$ws.Range("K20").Value = 3.75
$ws.Range("N20").Value = 3.25
$ws.Range("O20").Value = 1.4
$ws.Range("P20").Value = 1.77
$ws.Range("Q20").Value = 2.16
$ws.Range("S20").Value = 4.1
$ws.Range("W20").Value = 1.26
$ws.Range("Y20").Value = 9.4
# Row 21
$ws.Range("AN21").Value = 40
$ws.Range("AO21").Value = 30
$ws.Range("G21").Value = 3.05
$ws.Range("H21").Value = 2.58
$ws.Range("I21").Value = 2.66
$ws.Range("J21").Value = 3.35
$ws.Range("K21").Value = 3.5
$ws.Range("M21").Value = 1.09
$ws.Range("N21").Value = 3.25
$ws.Range("Q21").Value = 2.22
$ws.Range("R21").Value = 1.28
$ws.Range("T21").Value = 1.9
$ws.Range("U21").Value = 2.02
$ws.Range("V21").Value = 1.6
$ws.Range("W21").Value = 1.48
# Row 22
$ws.Range("F22").Value = 2.26
$ws.Range("G22").Value = 2.44
$ws.Range("I22").Value = 4
$ws.Range("K22").Value = 3.5
$ws.Range("N22").Value = 3
$ws.Range("Y22").Value = 14.5
